# Corrections to the vendor/origin lookup table:
#  - row 2 ("Erik" / SSA) becomes "Naiara" / SAO data with a Package 2 table
#    type, and the cnpj column switches from a free-text value to a real
#    numeric cnpj (plain "0" number format).
#  - row 3, previously blank, is filled in with a second matching record for
#    the same vendor/cnpj/origin but the ".COM 2" table type.
#  - a new right-aligned cell further out on the sheet (J11) extends the
#    used range, matching the refreshed selection/dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: fix the vendor/cnpj/table/origin data -------------------------
$ws.Range("A2").Value = "Naiara"

$ws.Range("B2").NumberFormat = "0"
$ws.Range("B2").Value = 30466928000127

$ws.Range("D2").Value = "SAO"
$ws.Range("C2").Value = "Package 2"

# --- Row 3: populate the previously empty row with the sibling record -----
$ws.Range("A3").Value = "Naiara"

$ws.Range("B3").NumberFormat = "0"
$ws.Range("B3").Value = 30466928000127

$ws.Range("D3").Value = "SAO"
$ws.Range("C3").Value = ".COM 2"
$ws.Range("E3").Value = "N"
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0

# --- New far-out formatted cell, extends the used range to J11 ------------
# -4152 == xlRight
$ws.Range("J11").HorizontalAlignment = -4152

# --- Column widths refreshed (best-fit) for the new content ----------------
$ws.Columns.Item(2).ColumnWidth = 14.3
$ws.Columns.Item(5).ColumnWidth = 11.17

# --- Selection moved to G7 -------------------------------------------------
$ws.Range("G7").Select() | Out-Null
